$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the tn (E) and fn (F) columns to 0 for rows 8, 9, 10
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
